$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseSteps")

# Delete rows 4-8, keeping only header + 2 data rows
$ws.Rows("4:8").Delete()

# Update remaining data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Login_Page"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Add_User"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1

# Update selection to match diff (single cell D9)
$ws.Range("D9").Select()
